# "udpated xlsx files issue n.34"
# Append 5 new coded-segment rows (315-319) to Sheet1, styled like the
# existing last row (314), and leave the active selection on D10 - matching
# the state the author's workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow  = 314
$startRow = 315
$endRow   = 319

$hospitalName = "San Gerardo de" + [char]0x2019 + " Tintori" + [char]0x201D + " Hospital"

# Column B (Comment) and C (Document group) are blank text cells in the
# source data; Document name (D) is a text-formatted number. A leading
# apostrophe forces Excel to store these as literal text (shared-string)
# values instead of numbers/blanks.
$newRows = @(
    @("●", "'", "'", "'18945", "Location:Country",       "1: 1660", "1: 1663", 0, "Iran",        4,  0.026444532592886424, "dattaray", "9/17/2019 11:38:55"),
    @("●", "'", "'", "'10622", "Location:Country",       "1: 222",  "1: 226",  0, "Italy",       5,  0.044385264092321346, "dattaray", "9/17/2019 11:39:49"),
    @("●", "'", "'", "'10622", "Location:City",          "1: 210",  "1: 214",  0, "Monza",       5,  0.044385264092321346, "dattaray", "9/17/2019 11:40:35"),
    @("●", "'", "'", "'10622", "Location:Hospital name", "1: 156",  "1: 188",  0, $hospitalName, 33, 0.29294274300932088,  "dattaray", "9/17/2019 11:40:57"),
    @("●", "'", "'", "'137",   "Location:Country",       "1: 454",  "1: 459",  0, "France",      6,  0.027466239414053559, "dattaray", "9/17/2019 11:41:56")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Copy the previous row's formatting (fill/border/font/number format/row
# height) down onto the new rows without disturbing the values just set.
$ws.Range("A$lastRow`:M$lastRow").Copy()
$ws.Range("A$startRow`:M$endRow").PasteSpecial(-4122)
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.6
}

# Restore the selection the author left active when they saved.
$ws.Range("D10").Select()

$wb.Save()
